# Update the cryptocurrency price/volume figures to the latest scrape,
# and swap the Maker/Cosmos rows (43/44) to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.907.51"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +5.76%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'3.664.03"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +17.34%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'596.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +3.01%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'182.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +4.57%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'3.660.93"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +17.38%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +0.12%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.535"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +3.63%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.162"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +6.31%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'6.61"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +3.20%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.497"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +4.65%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'40.24"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +11.28%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.0000253"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +4.67%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'4.287.19"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +17.67%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.675.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +17.95%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'70.967.70"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +5.90%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D19").Value = "'7.47"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +6.49%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'16.96"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.77%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'511.73"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +5.48%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'9.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +16.73%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.740"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +6.51%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'87.43"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +4.21%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'2.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +9.03%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'13.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +4.83%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'10.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +7.15%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -0.08%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +9.59%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'8.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +1.44%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'2.77"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +6.33%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +17.28%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'31.40"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +11.70%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +3.09%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  +0.08%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'6.08"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +8.29%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'1.01"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +6.85%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.345"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +10.49%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'2.15"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +8.80%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'51.03"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +3.54%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.127"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +3.27%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'45.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -6.00%  "
$ws.Range("E42").ClearFormats()
$ws.Range("B43").Value = "'Maker"
$ws.Range("B43").ClearFormats()
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C43").ClearFormats()
$ws.Range("D43").Value = "'3.129.73"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +11.40%  "
$ws.Range("E43").ClearFormats()
$ws.Range("B44").Value = "'Cosmos"
$ws.Range("B44").ClearFormats()
$ws.Range("C44").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C44").ClearFormats()
$ws.Range("D44").Value = "'8.80"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +5.58%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'413.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +10.55%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'2.77"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +3.56%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.0369"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +5.66%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'28.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +14.15%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'137.54"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +2.06%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D51").Value = "'2.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +11.38%  "
$ws.Range("E51").ClearFormats()
